# Apply the "datos" worksheet updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay plain text even if it looks like a number or a
    # date (e.g. "234", "2024-09-13"), then drop the leftover
    # quote-prefix/number-format style so the cell is indistinguishable from
    # a normal, never-specially-formatted text cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

function Set-EmptyTextValue($range) {
    # A real empty-string write clears the cell entirely in Excel's COM
    # model, so use the lone quote-prefix trick to leave behind an
    # empty (but present) text cell, then strip the leftover style.
    $range.Value = "'"
    $range.Style = "Normal"
}

# --- Row 2: unidad_numero (E2) becomes a true number instead of text "5" ---
$ws.Range("E2").Value = 5

# --- Row 5: several fields updated ---
$ws.Range("B5").Value = "miernesssss"
Set-TextValue $ws.Range("C5") "2024-09-13"
$ws.Range("E5").Value = 5
$ws.Range("G5").Value = "S6A"
$ws.Range("H5").Value = "Práctica"
$ws.Range("I5").Value = "asdf"
$ws.Range("J5").Value = "sadf"
$ws.Range("K5").Value = "asdf"

# --- Row 8: id bumped to 11, several fields updated ---
$ws.Range("A8").Value = 11
Set-TextValue $ws.Range("B8") "234"
$ws.Range("E8").Value = 23
$ws.Range("G8").Value = "M1A"
$ws.Range("H8").Value = "Explicación"
$ws.Range("I8").Value = "safd"
$ws.Range("J8").Value = "asdf"
$ws.Range("K8").Value = "asdfas"

# --- Row 9: replaced with a new Jueves entry ---
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Jueves"
$ws.Range("E9").Value = 5
$ws.Range("G9").Value = "M1A"
$ws.Range("H9").Value = "Explicación"
$ws.Range("I9").Value = "Jfjdjd"
$ws.Range("J9").Value = "Fjfjdjf"
$ws.Range("K9").Value = "Jfnfjjd"

# --- Row 10: brand new row ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "miernes"
Set-TextValue $ws.Range("C10") "2024-09-13"
Set-EmptyTextValue $ws.Range("D10")
$ws.Range("E10").Value = 5
Set-EmptyTextValue $ws.Range("F10")
$ws.Range("G10").Value = "S6A"
$ws.Range("H10").Value = "Práctica"
$ws.Range("I10").Value = "asdf"
$ws.Range("J10").Value = "sadf"
$ws.Range("K10").Value = "asdf"

# --- Row 11: brand new row ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "BNBBB"
Set-TextValue $ws.Range("C11") "2024-09-13"
Set-EmptyTextValue $ws.Range("D11")
Set-TextValue $ws.Range("E11") "3232"
Set-EmptyTextValue $ws.Range("F11")
$ws.Range("G11").Value = "TM"
$ws.Range("H11").Value = "Examen"
$ws.Range("I11").Value = "sadf"
$ws.Range("J11").Value = "asdf"
$ws.Range("K11").Value = "asdf"
